$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.195012493665558
$ws.Range("C2").Value = 0.327957343722403
$ws.Range("D2").Value = 0.03309288413318257
$ws.Range("F2").Value = 0.2359276067140001
$ws.Range("G2").Value = 0.1078855475638427
$ws.Range("H2").Value = 0.2908243665929504
$ws.Range("I2").Value = 0.1798627263767711
$ws.Range("O2").Value = 0.6885568647115008
# Row 3
$ws.Range("B3").Value = 1.043721088755092
$ws.Range("C3").Value = 0.2908277482814015
$ws.Range("D3").Value = 0.02887247580594732
$ws.Range("F3").Value = 0.2385816126413189
$ws.Range("G3").Value = 0.1107947158929079
$ws.Range("H3").Value = 0.2959752681811096
$ws.Range("I3").Value = 0.185723227919536
$ws.Range("O3").Value = 0.705287198323262
# Row 4
$ws.Range("B4").Value = 0.9504172635526515
$ws.Range("C4").Value = 0.2679376871411137
$ws.Range("D4").Value = 0.02627159485948027
$ws.Range("F4").Value = 0.2405130863643876
$ws.Range("G4").Value = 0.1128013793998974
$ws.Range("H4").Value = 0.299363065846908
$ws.Range("I4").Value = 0.1895597620483036
$ws.Range("O4").Value = 0.7164909371449255
# Row 5
$ws.Range("B5").Value = 0.9122951964980643
$ws.Range("C5").Value = 0.2585874278479423
$ws.Range("D5").Value = 0.02520939956293944
$ws.Range("F5").Value = 0.241375862683153
$ws.Range("G5").Value = 0.1136742289914281
$ws.Range("H5").Value = 0.3008001567177452
$ws.Range("I5").Value = 0.191182890128144
$ws.Range("O5").Value = 0.7212899080044579
# Row 6
$ws.Range("B6").Value = 0.9059591102536046
$ws.Range("C6").Value = 0.2570334981167264
$ws.Range("D6").Value = 0.02503288542973792
$ws.Range("F6").Value = 0.2415236904206779
$ws.Range("G6").Value = 0.1138224848919513
$ws.Range("H6").Value = 0.3010421970289769
$ws.Range("I6").Value = 0.1914560103617262
$ws.Range("O6").Value = 0.7221008436751291
# Row 7
$ws.Range("B7").Value = 0.9499035367047099
$ws.Range("C7").Value = 0.2678116757394946
$ws.Range("D7").Value = 0.02625727898713848
$ws.Range("F7").Value = 0.2405244159283981
$ws.Range("G7").Value = 0.1128129281946606
$ws.Range("H7").Value = 0.299382218143144
$ws.Range("I7").Value = 0.1895814105971727
$ws.Range("O7").Value = 0.716554713973224
# Row 8
$ws.Range("B8").Value = 1.142934222112103
$ws.Range("C8").Value = 0.3151747186313401
$ws.Range("D8").Value = 0.03163971378531016
$ws.Range("F8").Value = 0.2367799093135154
$ws.Range("G8").Value = 0.1088426950969676
$ws.Range("H8").Value = 0.2925536297221214
$ws.Range("I8").Value = 0.1818338901303411
$ws.Range("O8").Value = 0.6941317176824882
# Row 9
$ws.Range("B9").Value = 1.518090439092873
$ws.Range("C9").Value = 0.4072881913667743
$ws.Range("D9").Value = 0.04211590599184944
$ws.Range("F9").Value = 0.2318429293622373
$ws.Range("G9").Value = 0.1028197904265085
$ws.Range("H9").Value = 0.280952169741024
$ws.Range("I9").Value = 0.1685384926978779
$ws.Range("O9").Value = 0.6575852366608856
# Row 10
$ws.Range("B10").Value = 1.791523916020253
$ws.Range("C10").Value = 0.4744601906764387
$ws.Range("D10").Value = 0.04976131823518415
$ws.Range("F10").Value = 0.2296971834258414
$ws.Range("G10").Value = 0.0994881274531707
$ws.Range("H10").Value = 0.2735235100196789
$ws.Range("I10").Value = 0.1599374023737221
$ws.Range("O10").Value = 0.6353096062943706
# Row 11
$ws.Range("B11").Value = 1.91541280536012
$ws.Range("C11").Value = 0.5049016903568031
$ws.Range("D11").Value = 0.05322756170164666
$ws.Range("F11").Value = 0.2290458425457587
$ws.Range("G11").Value = 0.09821386278927235
$ws.Range("H11").Value = 0.270382633483905
$ws.Range("I11").Value = 0.156280140333366
$ws.Range("O11").Value = 0.6261795788326197
# Row 12
$ws.Range("B12").Value = 1.962251998559509
$ws.Range("C12").Value = 0.5164117504726278
$ws.Range("D12").Value = 0.0545383833739379
$ws.Range("F12").Value = 0.2288461437004941
$ws.Range("G12").Value = 0.09776635674497669
$ws.Range("H12").Value = 0.2692276321599536
$ws.Range("I12").Value = 0.1549321477180339
$ws.Range("O12").Value = 0.6228673931974384
# Row 13
$ws.Range("B13").Value = 1.952167721342619
$ws.Range("C13").Value = 0.5139336442889544
$ws.Range("D13").Value = 0.05425615457097877
$ws.Range("F13").Value = 0.2288870603905053
$ws.Range("G13").Value = 0.09786117199929834
$ws.Range("H13").Value = 0.2694748518472565
$ws.Range("I13").Value = 0.1552208161845838
$ws.Range("O13").Value = 0.6235742616331379
# Row 14
$ws.Range("B14").Value = 1.919267812801252
$ws.Range("C14").Value = 0.5058489851924151
$ws.Range("D14").Value = 0.05333543975606858
$ws.Range("F14").Value = 0.2290284708450727
$ws.Range("G14").Value = 0.09817634226764227
$ws.Range("H14").Value = 0.2702869212785259
$ws.Range("I14").Value = 0.1561684988049166
$ws.Range("O14").Value = 0.6259041690850751
# Row 15
$ws.Range("B15").Value = 1.899105824388585
$ws.Range("C15").Value = 0.5008945961739073
$ws.Range("D15").Value = 0.05277124229816366
$ws.Range("F15").Value = 0.2291212104047773
$ws.Range("G15").Value = 0.09837396483197125
$ws.Range("H15").Value = 0.2707888172766602
$ws.Range("I15").Value = 0.1567537976193609
$ws.Range("O15").Value = 0.6273502347441138
# Row 16
$ws.Range("B16").Value = 1.783417192997092
$ws.Range("C16").Value = 0.4724683696379657
$ws.Range("D16").Value = 0.04953454824359937
$ws.Range("F16").Value = 0.229746306572288
$ws.Range("G16").Value = 0.09957628782971284
$ws.Range("H16").Value = 0.2737335798226184
$ws.Range("I16").Value = 0.1601815701853946
$ws.Range("O16").Value = 0.6359265395235667
# Row 17
$ws.Range("B17").Value = 1.712316205624575
$ws.Range("C17").Value = 0.4549996177792082
$ws.Range("D17").Value = 0.04754588365759105
$ws.Range("F17").Value = 0.2302131417114879
$ws.Range("G17").Value = 0.10037591749564
$ws.Range("H17").Value = 0.2756012356755235
$ws.Range("I17").Value = 0.162349961727621
$ws.Range("O17").Value = 0.6414454204320208
# Row 18
$ws.Range("B18").Value = 1.671374179944735
$ws.Range("C18").Value = 0.4449412460112399
$ws.Range("D18").Value = 0.04640096103616997
$ws.Range("F18").Value = 0.2305122071320156
$ws.Range("G18").Value = 0.1008585385842125
$ws.Range("H18").Value = 0.2766979013562647
$ws.Range("I18").Value = 0.1636211920646389
$ws.Range("O18").Value = 0.6447141144413138
# Row 19
$ws.Range("B19").Value = 1.657504008393857
$ws.Range("C19").Value = 0.441533824964722
$ws.Range("D19").Value = 0.04601312447077532
$ws.Range("F19").Value = 0.2306187056451279
$ws.Range("G19").Value = 0.101025833688162
$ws.Range("H19").Value = 0.2770730647387083
$ws.Range("I19").Value = 0.1640557295574077
$ws.Range("O19").Value = 0.6458370206766233
# Row 20
$ws.Range("B20").Value = 1.719889876057721
$ws.Range("C20").Value = 0.4568603212208018
$ws.Range("D20").Value = 0.04775769440226441
$ws.Range("F20").Value = 0.2301602822621049
$ws.Range("G20").Value = 0.1002884442954439
$ws.Range("H20").Value = 0.27540009746739
$ws.Range("I20").Value = 0.1621166443653008
$ws.Range("O20").Value = 0.6408481517358382
# Row 21
$ws.Range("B21").Value = 1.928933370906407
$ws.Range("C21").Value = 0.5082241268207781
$ws.Range("D21").Value = 0.05360592474002601
$ws.Range("F21").Value = 0.2289856589973098
$ws.Range("G21").Value = 0.09808281580544786
$ws.Range("H21").Value = 0.2700474630084528
$ws.Range("I21").Value = 0.1558891374652012
$ws.Range("O21").Value = 0.6252158722030572
# Row 22
$ws.Range("B22").Value = 2.065117631210001
$ws.Range("C22").Value = 0.5416910795636909
$ws.Range("D22").Value = 0.05741772671812839
$ws.Range("F22").Value = 0.228491747461419
$ws.Range("G22").Value = 0.09684564406059337
$ws.Range("H22").Value = 0.2667496440779402
$ws.Range("I22").Value = 0.152034470861909
$ws.Range("O22").Value = 0.615845771165084
# Row 23
$ws.Range("B23").Value = 1.992474671905768
$ws.Range("C23").Value = 0.5238387842398993
$ws.Range("D23").Value = 0.05538427271886803
$ws.Range("F23").Value = 0.228730225804668
$ws.Range("G23").Value = 0.09748714161343486
$ws.Range("H23").Value = 0.2684913823744495
$ws.Range("I23").Value = 0.1540720052553501
$ws.Range("O23").Value = 0.620769023966119
# Row 24
$ws.Range("B24").Value = 1.716466022245697
$ws.Range("C24").Value = 0.4560191449538138
$ws.Range("D24").Value = 0.04766193978210254
$ws.Range("F24").Value = 0.2301840844868863
$ws.Range("G24").Value = 0.1003279196126954
$ws.Range("H24").Value = 0.2754909606472395
$ws.Range("I24").Value = 0.1622220505550578
$ws.Range("O24").Value = 0.6411178785243976
# Row 25
$ws.Range("B25").Value = 1.416976320814911
$ws.Range("C25").Value = 0.3824548638710326
$ws.Range("D25").Value = 0.03929060586338551
$ws.Range("F25").Value = 0.2329194067134601
$ws.Range("G25").Value = 0.1042584865482041
$ws.Range("H25").Value = 0.2838986827935237
$ws.Range("I25").Value = 0.1719310481100571
$ws.Range("O25").Value = 0.6666720873480827
